$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# List of cell updates produced by the "Update automatic: dades i banners [2026-02-18 20:50]" run.
# Each entry sets a cell's content as text (via .Formula). Values that would otherwise be
# auto-interpreted by Excel as a pure number (e.g. "73%") are prefixed with a leading
# apostrophe so they are stored as literal text, matching the source data (all cells in this
# sheet are plain text, not real Excel numbers/percentages/dates).
$updates = @(
    @{ Cell = 'E2'; Value = '2026-02-18 20:48:20' },
    @{ Cell = 'E3'; Value = '2026-02-18 20:48:23' },
    @{ Cell = 'E4'; Value = '2026-02-18 20:48:25' },
    @{ Cell = 'H4'; Value = '''73%' },
    @{ Cell = 'J4'; Value = '1013.0 hPa' },
    @{ Cell = 'E5'; Value = '2026-02-18 20:48:27' },
    @{ Cell = 'H5'; Value = '''66%' },
    @{ Cell = 'O5'; Value = '0.7 °C' },
    @{ Cell = 'E6'; Value = '2026-02-18 20:48:30' },
    @{ Cell = 'J6'; Value = '1012.7 hPa' },
    @{ Cell = 'O6'; Value = '12.1 °C' },
    @{ Cell = 'E7'; Value = '2026-02-18 20:48:32' },
    @{ Cell = 'J7'; Value = '1014.2 hPa' },
    @{ Cell = 'E8'; Value = '2026-02-18 20:48:35' },
    @{ Cell = 'J8'; Value = '1014.0 hPa' },
    @{ Cell = 'E9'; Value = '2026-02-18 20:48:37' },
    @{ Cell = 'H9'; Value = '''78%' },
    @{ Cell = 'O9'; Value = '11.0 °C' },
    @{ Cell = 'E10'; Value = '2026-02-18 20:48:40' },
    @{ Cell = 'K10'; Value = '12.0 MJ/m2' },
    @{ Cell = 'E11'; Value = '2026-02-18 20:48:42' },
    @{ Cell = 'O11'; Value = '5.7 °C' },
    @{ Cell = 'E12'; Value = '2026-02-18 20:48:45' },
    @{ Cell = 'H12'; Value = '''87%' },
    @{ Cell = 'O12'; Value = '10.9 °C' },
    @{ Cell = 'E13'; Value = '2026-02-18 20:48:47' },
    @{ Cell = 'J13'; Value = '1015.3 hPa' },
    @{ Cell = 'O13'; Value = '3.9 °C' },
    @{ Cell = 'E14'; Value = '2026-02-18 20:48:49' },
    @{ Cell = 'H14'; Value = '''87%' },
    @{ Cell = 'E15'; Value = '2026-02-18 20:48:52' },
    @{ Cell = 'E16'; Value = '2026-02-18 20:48:54' },
    @{ Cell = 'H16'; Value = '''52%' },
    @{ Cell = 'E17'; Value = '2026-02-18 20:48:57' },
    @{ Cell = 'O17'; Value = '3.5 °C' },
    @{ Cell = 'E18'; Value = '2026-02-18 20:48:59' },
    @{ Cell = 'J18'; Value = '1013.2 hPa' },
    @{ Cell = 'E19'; Value = '2026-02-18 20:49:02' },
    @{ Cell = 'K19'; Value = '9.8 MJ/m2' },
    @{ Cell = 'E20'; Value = '2026-02-18 20:49:04' },
    @{ Cell = 'H20'; Value = '''75%' },
    @{ Cell = 'O20'; Value = '-0.4 °C' },
    @{ Cell = 'E21'; Value = '2026-02-18 20:49:07' },
    @{ Cell = 'J21'; Value = '1014.8 hPa' },
    @{ Cell = 'E22'; Value = '2026-02-18 20:49:09' },
    @{ Cell = 'I22'; Value = '0.6 mm' },
    @{ Cell = 'E23'; Value = '2026-02-18 20:49:12' },
    @{ Cell = 'H23'; Value = '''56%' },
    @{ Cell = 'N23'; Value = '-3.0 °C 20:25 TU' },
    @{ Cell = 'E24'; Value = '2026-02-18 20:49:14' },
    @{ Cell = 'H24'; Value = '''86%' },
    @{ Cell = 'J24'; Value = '1014.8 hPa' },
    @{ Cell = 'E25'; Value = '2026-02-18 20:49:17' },
    @{ Cell = 'E26'; Value = '2026-02-18 20:49:19' },
    @{ Cell = 'G26'; Value = '1 cm' },
    @{ Cell = 'J26'; Value = '1012.0 hPa' },
    @{ Cell = 'O26'; Value = '5.6 °C' },
    @{ Cell = 'E27'; Value = '2026-02-18 20:49:22' },
    @{ Cell = 'H27'; Value = '''56%' },
    @{ Cell = 'E28'; Value = '2026-02-18 20:49:24' },
    @{ Cell = 'J28'; Value = '1012.8 hPa' },
    @{ Cell = 'O28'; Value = '9.9 °C' },
    @{ Cell = 'E29'; Value = '2026-02-18 20:49:27' },
    @{ Cell = 'O29'; Value = '12.0 °C' },
    @{ Cell = 'E30'; Value = '2026-02-18 20:49:29' },
    @{ Cell = 'H30'; Value = '''76%' },
    @{ Cell = 'J30'; Value = '1012.4 hPa' },
    @{ Cell = 'O30'; Value = '11.1 °C' },
    @{ Cell = 'E31'; Value = '2026-02-18 20:49:32' },
    @{ Cell = 'J31'; Value = '1011.1 hPa' },
    @{ Cell = 'E32'; Value = '2026-02-18 20:49:34' },
    @{ Cell = 'H32'; Value = '''84%' },
    @{ Cell = 'E33'; Value = '2026-02-18 20:49:37' },
    @{ Cell = 'J33'; Value = '1014.0 hPa' },
    @{ Cell = 'E34'; Value = '2026-02-18 20:49:39' },
    @{ Cell = 'H34'; Value = '''47%' },
    @{ Cell = 'E35'; Value = '2026-02-18 20:49:42' },
    @{ Cell = 'J35'; Value = '1014.3 hPa' },
    @{ Cell = 'O35'; Value = '9.3 °C' },
    @{ Cell = 'E36'; Value = '2026-02-18 20:49:44' },
    @{ Cell = 'J36'; Value = '1012.9 hPa' },
    @{ Cell = 'O36'; Value = '11.8 °C' },
    @{ Cell = 'E37'; Value = '2026-02-18 20:49:47' },
    @{ Cell = 'J37'; Value = '1014.5 hPa' },
    @{ Cell = 'E38'; Value = '2026-02-18 20:49:49' },
    @{ Cell = 'E39'; Value = '2026-02-18 20:49:52' },
    @{ Cell = 'H39'; Value = '''41%' },
    @{ Cell = 'E40'; Value = '2026-02-18 20:49:54' },
    @{ Cell = 'J40'; Value = '1015.4 hPa' },
    @{ Cell = 'O40'; Value = '6.3 °C' },
    @{ Cell = 'E41'; Value = '2026-02-18 20:49:56' },
    @{ Cell = 'J41'; Value = '1014.5 hPa' },
    @{ Cell = 'O41'; Value = '11.4 °C' },
    @{ Cell = 'E42'; Value = '2026-02-18 20:49:59' },
    @{ Cell = 'O42'; Value = '11.8 °C' },
    @{ Cell = 'E43'; Value = '2026-02-18 20:50:01' },
    @{ Cell = 'E44'; Value = '2026-02-18 20:50:04' },
    @{ Cell = 'H44'; Value = '''73%' },
    @{ Cell = 'E45'; Value = '2026-02-18 20:50:06' },
    @{ Cell = 'H45'; Value = '''63%' },
    @{ Cell = 'J45'; Value = '1011.7 hPa' },
    @{ Cell = 'O45'; Value = '7.5 °C' },
    @{ Cell = 'E46'; Value = '2026-02-18 20:50:08' },
    @{ Cell = 'J46'; Value = '1014.8 hPa' }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Formula = $u.Value
}
